# Cable price change: zero out the INVCOST / FIXOM cable trade values
# (formerly small non-zero placeholders) on the ELC_TRADE sheet.
# Dependent formula cells (J/K/L/M columns) recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ELC_TRADE")

# TB_ELCC_DKE_DKISLBH_01 block
$ws.Range("H11").Value = 0   # INVCOST
$ws.Range("H12").Value = 0   # FIXOM

# TB_ELCC_DKW_DKISL1_01 block
$ws.Range("I16").Value = 0   # INVCOST
$ws.Range("I17").Value = 0   # FIXOM

# TB_ELCC_DKW_DKISL2_01 block
$ws.Range("I21").Value = 0   # INVCOST
$ws.Range("I22").Value = 0   # FIXOM

# TB_ELCC_DKW_DKISL3_01 block
$ws.Range("I26").Value = 0   # INVCOST
$ws.Range("I27").Value = 0   # FIXOM

$excel.Calculate()
